# Assignment 3 Iteration: Generics
# Rebuild the JGRAM RESULTS table: merge the stray duplicate <w:tblPr>,
# mark the table layout as fixed, give every cell an explicit <w:tcW>
# (score columns 1000 twips wide, the feedback column 6000 twips wide),
# and refresh the embedded grading token in the last cell.

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$tableRange = $table.Range
$insertAt = $tableRange.Start

$table.Delete()

$newTableXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblW w:w="0" w:type="auto"/><w:tblBorders><w:top w:val="single"/><w:left w:val="single"/><w:bottom w:val="single"/><w:right w:val="single"/><w:insideH w:val="single"/><w:insideV w:val="single"/></w:tblBorders><w:tblLayout w:type="fixed"/></w:tblPr><w:tr><w:tc><w:tcPr><w:shd w:color="auto" w:val="clear" w:fill="c0c0c0"/></w:tcPr><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>C#</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:shd w:color="auto" w:val="clear" w:fill="c0c0c0"/></w:tcPr><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>Weight</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:shd w:color="auto" w:val="clear" w:fill="c0c0c0"/></w:tcPr><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>Grade</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:shd w:color="auto" w:val="clear" w:fill="c0c0c0"/></w:tcPr><w:tcPr><w:tcW w:w="6000"/></w:tcPr><w:p><w:r><w:t>Feedback</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>85</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6000"/></w:tcPr><w:p><w:r><w:t>Throws ArrayIndexOutOfBoundsException; watch out for the Boolean condition that controls the for loop’s execution. This for loop executes one more time than you would want it to because of the greater than or equal to sign.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>100</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6000"/></w:tcPr><w:p><w:r><w:t>Excellent work. Please make sure to include comments next time.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>90</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6000"/></w:tcPr><w:p><w:r><w:t>Per the prompt, the method should return the first element in the array; use the break keyword to exit the for loop once the element is found.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t/></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>Σ</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:shd w:color="auto" w:val="clear" w:fill="8fbc8f"/></w:tcPr><w:tcPr><w:tcW w:w="1000"/></w:tcPr><w:p><w:r><w:t>91.50</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6000"/></w:tcPr><w:p><w:r><w:t>eyJhbGciOiJIUzI1NiJ9.eyJqdGkiOiIxIiwiaWF0IjoxNTkwOTQwMjM0LCJzdWIiOiJKR1JBTSIsImlzcyI6IkJVLU1FVCIsIjEtV2VpZ2h0IjozLCIxLUdyYWRlIjo4NSwiMS1GZWVkYmFjayI6IlRocm93cyBBcnJheUluZGV4T3V0T2ZCb3VuZHNFeGNlcHRpb247IHdhdGNoIG91dCBmb3IgdGhlIEJvb2xlYW4gY29uZGl0aW9uIHRoYXQgY29udHJvbHMgdGhlIGZvciBsb29w4oCZcyBleGVjdXRpb24uIFRoaXMgZm9yIGxvb3AgZXhlY3V0ZXMgb25lIG1vcmUgdGltZSB0aGFuIHlvdSB3b3VsZCB3YW50IGl0IHRvIGJlY2F1c2Ugb2YgdGhlIGdyZWF0ZXIgdGhhbiBvciBlcXVhbCB0byBzaWduLiIsIjItV2VpZ2h0IjozLCIyLUdyYWRlIjoxMDAsIjItRmVlZGJhY2siOiJFeGNlbGxlbnQgd29yay4gUGxlYXNlIG1ha2Ugc3VyZSB0byBpbmNsdWRlIGNvbW1lbnRzIG5leHQgdGltZS4iLCIzLVdlaWdodCI6NCwiMy1HcmFkZSI6OTAsIjMtRmVlZGJhY2siOiJQZXIgdGhlIHByb21wdCwgdGhlIG1ldGhvZCBzaG91bGQgcmV0dXJuIHRoZSBmaXJzdCBlbGVtZW50IGluIHRoZSBhcnJheTsgdXNlIHRoZSBicmVhayBrZXl3b3JkIHRvIGV4aXQgdGhlIGZvciBsb29wIG9uY2UgdGhlIGVsZW1lbnQgaXMgZm91bmQuIiwiQ1BJbmRleGVzIjoiWzEsIDIsIDNdIiwiR3JhZGVNYXBwaW5nIjoiQSsgPSA5N1xuQSAgPSA5NVxuQS0gPSA5M1xuQisgPSA4N1xuQiAgPSA4NVxuQi0gPSA4M1xuQysgPSA3N1xuQyAgPSA3NVxuQy0gPSA3M1xuRiAgPSA2N1xuIiwiVG90YWxHcmFkZSI6OTEuNX0.3U49MndlEDPMI6GjPD6gN_sKLQQOrgv_xlA4P88CXxc</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target = $d.Range($insertAt, $insertAt)
$target.InsertXML($newTableXml)

Write-Host "Table rebuilt; table count:" $d.Tables.Count
